$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '292.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-6.63%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.43'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.90%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.006'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.96%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07323'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.39%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.545'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-8.12%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9245'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.09%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.384'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.65%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1209'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.74%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1738'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-4.35%'

$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04309'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.43%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08596'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.16%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1055'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.12%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001273'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.88%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005874'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.58%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.338'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.33%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.292'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.92%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.693'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.72%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.07%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2794'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.72%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.03929'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.01%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001260'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.73%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003782'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-6.70%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001281'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.77%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003725'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-95.05%'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02297'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-5.66%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.04975'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.86%'

$ws.Range("B40").Value = 'CEJI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.005374'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '85.53%'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007692'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.78%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1286'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.13%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007330'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.81%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007908'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.62%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3165'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.66%'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006273'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.70%'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.02045'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-92.41%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.01%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
